# This workbook contains a weekly price table for "Achicoria" at
# "Vega Modelo de Temuco". A new weekly record needs to be inserted as the
# new row 5 (immediately after the header-adjacent rows 2-4), pushing all
# the existing price rows (previously rows 5-26) down by one row, and
# finishing with the used range growing from A1:R26 to A1:R27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 5; this shifts rows 5:26 down to
# 6:27 and carries formatting (incl. the date style on column D) along.
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new weekly record.
$ws.Range("A5").Value = 10
$ws.Range("B5").Value = "Vega Modelo de Temuco"
$ws.Range("C5").Value = "La Araucanía"
$ws.Range("D5").Value = 44676
$ws.Range("E5").Value = 9
$ws.Range("F5").Value = 100112010
$ws.Range("G5").Value = "Achicoria"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 40
$ws.Range("K5").Value = 12000
$ws.Range("L5").Value = 12000
$ws.Range("M5").Value = 12000
$ws.Range("N5").Value = "$/caja 18 unidades"
$ws.Range("O5").Value = "Región Metropolitana"
$ws.Range("P5").Value = 667
$ws.Range("Q5").Value = 18
$ws.Range("R5").Value = "Hortaliza"
